{"js": "// \"deleted one page break\" \u2014 the paragraph mark that separated the two\n// trailing empty \"Heading 1\" placeholder paragraphs (the ones that used to\n// force page breaks before the \"\u0426\u0435\u043b\u044c \u043b\u0430\u0431\u043e\u0440\u0430\u0442\u043e\u0440\u043d\u043e\u0439 \u0440\u0430\u0431\u043e\u0442\u044b\" / \"\u0417\u0430\u0434\u0430\u043d\u0438\u0435\n// \u043b\u0430\u0431\u043e\u0440\u0430\u0442\u043e\u0440\u043d\u043e\u0439 \u0440\u0430\u0431\u043e\u0442\u044b\" headings) was removed. Word carried the now-merged\n// paragraph mark's direct formatting (left alignment / zero indent / en-US\n// language) onto both of the remaining empty Heading-1 paragraphs at the\n// end of the document.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/style,items/text\");\nawait context.sync();\n\n// Identify the trailing empty \"Heading 1\" paragraphs (there are exactly two\n// of them at the tail of the document, one of which carries the\n// section break). Walk from the end so the logic is resilient to any\n// unrelated paragraphs elsewhere in the body.\nconst items = paragraphs.items;\nconst targets = [];\nfor (let i = items.length - 1; i >= 0 && targets.length < 2; i--) {\n  const p = items[i];\n  if (p.style === \"Heading 1\" && p.text.trim() === \"\") {\n    targets.unshift(p);\n  } else if (targets.length > 0) {\n    // Stop once we've found the trailing run and hit something else.\n    break;\n  }\n}\n\ntargets.forEach((p, i) => {\n  p.alignment = Word.Alignment.left;\n  if (i === targets.length - 1) {\n    // Only the very last paragraph in the body also picks up the explicit\n    // zero left indent.\n    p.leftIndent = 0;\n  }\n  const rng = p.getRange();\n  rng.languageId = \"en-US\";\n});\n\nawait context.sync();\n", "ps1": "# \"deleted one page break\" \u2014 the paragraph mark that separated the two\n# trailing empty \"Heading 1\" placeholder paragraphs (the ones that used to\n# force page breaks before the \"\u0426\u0435\u043b\u044c \u043b\u0430\u0431\u043e\u0440\u0430\u0442\u043e\u0440\u043d\u043e\u0439 \u0440\u0430\u0431\u043e\u0442\u044b\" / \"\u0417\u0430\u0434\u0430\u043d\u0438\u0435\n# \u043b\u0430\u0431\u043e\u0440\u0430\u0442\u043e\u0440\u043d\u043e\u0439 \u0440\u0430\u0431\u043e\u0442\u044b\" headings) was removed. Word carried the now-merged\n# paragraph mark's direct formatting (left alignment / zero indent / en-US\n# language) onto both of the remaining empty Heading-1 paragraphs at the\n# end of the document.\n\n$d = $word.ActiveDocument\n\n# Identify the trailing empty \"Heading 1\" paragraphs (there are exactly two\n# of them at the tail of the document, one of which carries the\n# section break). Walk from the end so the logic is resilient to any\n# unrelated paragraphs elsewhere in the body.\n$targets = New-Object System.Collections.ArrayList\nfor ($i = $d.Paragraphs.Count; $i -ge 1; $i--) {\n    $p = $d.Paragraphs.Item($i)\n    $isEmptyHeading = ($p.Style.NameLocal -eq \"Heading 1\") -and ($p.Range.Text.Trim() -eq \"\")\n    if ($isEmptyHeading) {\n        [void]$targets.Insert(0, $p)\n        if ($targets.Count -eq 2) { break }\n    } elseif ($targets.Count -gt 0) {\n        break\n    }\n}\n\nfor ($i = 0; $i -lt $targets.Count; $i++) {\n    $p = $targets[$i]\n    $p.Alignment = \"wdAlignParagraphLeft\"\n    if ($i -eq $targets.Count - 1) {\n        # Only the very last paragraph in the body also picks up the\n        # explicit zero left indent.\n        $p.LeftIndent = 0\n    }\n    $p.Range.LanguageID = \"en-US\"\n}\n"}
